# Auto-generated edit script applying the Masamune_Profits.xlsx data refresh
# Updates numeric profit/cost columns (H, I, J, K, L, M, N) on several sheets
# to match the latest scheduled-runner snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 50593
$ws.Range("J93").Value = 50593
$ws.Range("L93").Value = 50593
$ws.Range("N93").Value = -55585
$ws.Range("H95").Value = 35739.332
$ws.Range("J95").Value = 35739.332
$ws.Range("L95").Value = 35739.332
$ws.Range("N95").Value = -41231.332
$ws.Range("H105").Value = 49663
$ws.Range("J105").Value = 49663
$ws.Range("L105").Value = 49663
$ws.Range("N105").Value = -56651
$ws.Range("H107").Value = 1546.2142
$ws.Range("I107").Value = 430.625
$ws.Range("J107").Value = 3033.6667
$ws.Range("K107").Value = 430.625
$ws.Range("L107").Value = 3033.6667
$ws.Range("M107").Value = 1489.375
$ws.Range("N107").Value = -6873.6667
$ws.Range("H138").Value = 3667.3
$ws.Range("I138").Value = 2985.35
$ws.Range("J138").Value = 3940.08
$ws.Range("K138").Value = 8956.049999999999
$ws.Range("L138").Value = 11820.24
$ws.Range("M138").Value = -3816.049999999999
$ws.Range("N138").Value = -22100.24

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 49598
$ws.Range("J101").Value = 49598
$ws.Range("L101").Value = 49598
$ws.Range("N101").Value = -56088
$ws.Range("H103").Value = 38354
$ws.Range("J103").Value = 38354
$ws.Range("L103").Value = 38354
$ws.Range("N103").Value = -40698
$ws.Range("H105").Value = 49370
$ws.Range("J105").Value = 49370
$ws.Range("L105").Value = 49370
$ws.Range("N105").Value = -56358
$ws.Range("H122").Value = 1618.1818
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 558.5714
$ws.Range("I22").Value = 749.3333
$ws.Range("J22").Value = 415.5
$ws.Range("K22").Value = 749.3333
$ws.Range("L22").Value = 415.5
$ws.Range("M22").Value = -399.3333
$ws.Range("N22").Value = -1115.5
$ws.Range("H31").Value = 5906.931
$ws.Range("I31").Value = 2417.5908
$ws.Range("J31").Value = 16873.428
$ws.Range("K31").Value = 2417.5908
$ws.Range("L31").Value = 16873.428
$ws.Range("M31").Value = -2122.5908
$ws.Range("N31").Value = -17463.428
$ws.Range("H34").Value = 5906.931
$ws.Range("I34").Value = 2417.5908
$ws.Range("J34").Value = 16873.428
$ws.Range("K34").Value = 2417.5908
$ws.Range("L34").Value = 16873.428
$ws.Range("M34").Value = -2215.5908
$ws.Range("N34").Value = -17277.428
$ws.Range("H41").Value = 36971.8
$ws.Range("J41").Value = 44464.75
$ws.Range("L41").Value = 44464.75
$ws.Range("N41").Value = -45320.75
$ws.Range("H50").Value = 41949.75
$ws.Range("J50").Value = 41949.75
$ws.Range("L50").Value = 41949.75
$ws.Range("N50").Value = -43199.75
$ws.Range("H92").Value = 44601
$ws.Range("J92").Value = 44601
$ws.Range("L92").Value = 44601
$ws.Range("N92").Value = -49593
$ws.Range("H96").Value = 80156
$ws.Range("J96").Value = 80156
$ws.Range("L96").Value = 80156
$ws.Range("N96").Value = -85648

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2004
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9691.5
$ws.Range("J46").Value = 10640
$ws.Range("L46").Value = 10640
$ws.Range("N46").Value = -10952
$ws.Range("H57").Value = 35151.668
$ws.Range("J57").Value = 36620
$ws.Range("L57").Value = 36620
$ws.Range("N57").Value = -38260
$ws.Range("H80").Value = 281589.22
$ws.Range("J80").Value = 3073.2727
$ws.Range("L80").Value = 3073.2727
$ws.Range("N80").Value = -5069.2727
$ws.Range("H83").Value = 281589.22
$ws.Range("J83").Value = 3073.2727
$ws.Range("L83").Value = 15366.3635
$ws.Range("N83").Value = -25350.3635
$ws.Range("H86").Value = 27139
$ws.Range("J86").Value = 27139
$ws.Range("L86").Value = 27139
$ws.Range("N86").Value = -29511
$ws.Range("H89").Value = 27139
$ws.Range("J89").Value = 27139
$ws.Range("L89").Value = 81417
$ws.Range("N89").Value = -93273
$ws.Range("H101").Value = 50657
$ws.Range("J101").Value = 50657
$ws.Range("L101").Value = 50657
$ws.Range("N101").Value = -57147
$ws.Range("H122").Value = 1966
$ws.Range("I122").Value = 1844
$ws.Range("J122").Value = 2108.3333
$ws.Range("K122").Value = 5532
$ws.Range("L122").Value = 6324.999899999999
$ws.Range("M122").Value = -3082
$ws.Range("N122").Value = -11224.9999
$ws.Range("H124").Value = 39390
$ws.Range("J124").Value = 39390
$ws.Range("L124").Value = 39390
$ws.Range("N124").Value = -49210
$ws.Range("H125").Value = 44322
$ws.Range("J125").Value = 44322
$ws.Range("L125").Value = 44322
$ws.Range("N125").Value = -49242
$ws.Range("H126").Value = 11789.909
$ws.Range("I126").Value = 12798.9
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 38396.7
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -35926.7
$ws.Range("N126").Value = -10040
$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 57504
$ws.Range("J12").Value = 57504
$ws.Range("L12").Value = 57504
$ws.Range("N12").Value = -57844
$ws.Range("H35").Value = 38333.332
$ws.Range("J35").Value = 38333.332
$ws.Range("L35").Value = 38333.332
$ws.Range("N35").Value = -39005.332
$ws.Range("H55").Value = 449.2143
$ws.Range("I55").Value = 382.41666
$ws.Range("J55").Value = 850
$ws.Range("K55").Value = 382.41666
$ws.Range("L55").Value = 850
$ws.Range("M55").Value = -209.41666
$ws.Range("N55").Value = -1196
$ws.Range("H68").Value = 2554.5557
$ws.Range("I68").Value = 2291.1667
$ws.Range("J68").Value = 3081.3333
$ws.Range("K68").Value = 2291.1667
$ws.Range("L68").Value = 3081.3333
$ws.Range("M68").Value = -1542.1667
$ws.Range("N68").Value = -4579.3333
$ws.Range("H71").Value = 2554.5557
$ws.Range("I71").Value = 2291.1667
$ws.Range("J71").Value = 3081.3333
$ws.Range("K71").Value = 11455.8335
$ws.Range("L71").Value = 15406.6665
$ws.Range("M71").Value = -7711.833500000001
$ws.Range("N71").Value = -22894.6665
$ws.Range("H105").Value = 49871.668
$ws.Range("J105").Value = 49871.668
$ws.Range("L105").Value = 49871.668
$ws.Range("N105").Value = -56859.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 85740.2
$ws.Range("J92").Value = 85740.2
$ws.Range("L92").Value = 85740.2
$ws.Range("N92").Value = -90732.2
$ws.Range("H95").Value = 35384
$ws.Range("J95").Value = 35384
$ws.Range("L95").Value = 35384
$ws.Range("N95").Value = -40876
$ws.Range("H98").Value = 41585
$ws.Range("J98").Value = 41585
$ws.Range("L98").Value = 41585
$ws.Range("N98").Value = -47575
